# Add a new "2022-Q3" quarterly sheet to the workbook and update the
# "总计" (summary) sheet with the corresponding new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by duplicating the "2022-Q2"
#    sheet (this keeps identical column layout / styles / widths) and
#    placing the copy immediately before it, then renaming it.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)

$q2After = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Item($q2After.Index - 1)
$q3.Name = "2022-Q3"

# Fund rows for the 2022-Q3 sheet (code / name / size / position% /
# weight / value / rank). Columns B-G are stored as text in the
# original workbook, so force a text number-format before assigning
# the values to avoid Excel's automatic numeric coercion (which would
# strip leading zeros from fund codes, etc.).
$q3Rows = @(
    @{ Row = 2; Code = "210009"; Name = "金鹰核心资源混合";     Size = "2.84"; Pos = "93.42"; Pct = "5.16"; Val = "0.1465"; Rank = 8 },
    @{ Row = 3; Code = "162102"; Name = "金鹰中小盘精选混合"; Size = "3.17"; Pos = "78.28"; Pct = "4.61"; Val = "0.1461"; Rank = 6 },
    @{ Row = 4; Code = "001167"; Name = "金鹰科技创新股票";   Size = "2.66"; Pos = "94.84"; Pct = "5.12"; Val = "0.1362"; Rank = 9 }
)

foreach ($entry in $q3Rows) {
    $r = $entry.Row

    $cellB = $q3.Range("B$r")
    $cellB.NumberFormat = "@"
    $cellB.Value = $entry.Code

    $cellC = $q3.Range("C$r")
    $cellC.NumberFormat = "@"
    $cellC.Value = $entry.Name

    $cellD = $q3.Range("D$r")
    $cellD.NumberFormat = "@"
    $cellD.Value = $entry.Size

    $cellE = $q3.Range("E$r")
    $cellE.NumberFormat = "@"
    $cellE.Value = $entry.Pos

    $cellF = $q3.Range("F$r")
    $cellF.NumberFormat = "@"
    $cellF.Value = $entry.Pct

    $cellG = $q3.Range("G$r")
    $cellG.NumberFormat = "@"
    $cellG.Value = $entry.Val

    $q3.Range("H$r").Value = $entry.Rank
}

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: a new row is inserted for
#    2022-Q3 and every following row shifts down by one, adding a
#    final row for 2020-Q4.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Duplicate row 5's formatting into row 6 first (new last row), so that
# the A6 index cell keeps the same style ("s=2") as the other index
# cells in column A.
$total.Range("A5:D5").Copy($total.Range("A6:D6"))

# Now rewrite the label/value rows from the bottom up so each row
# simply reflects the row that used to be one above it, then put the
# brand-new 2022-Q3 summary into row 2.
$total.Range("B6").Value = "2020-Q4"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.02

$total.Range("B5").Value = "2021-Q1"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.07000000000000001

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.34

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.4

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.43

# Refresh the A-column running index (0..4) for all five rows.
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4

# Restore the originally-active sheet (the copy operation above shifts
# Excel's "active sheet" onto the newly created tab).
$wb.Worksheets.Item("2020-Q4").Activate()
